## Refactor invoice generation to use Manjerico data as the sender.
##
## "Dados Manjerico" (sheet2) held the sender's address as a single combined
## string ("Rua Itacuruçá, 26") in the "Endereço" column while the house
## number ("26") already lived in its own "Número" column. Split the street
## name away from the number (the number column is untouched) and add an
## "E-mail" column with Manjerico's contact address as a live mailto
## hyperlink, matching the columns already used for the other companies'
## records.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados Manjerico")

# Endereço no longer repeats the house number (Número/D2 already has it).
$ws.Range("C2").Value = "Rua Itacuruçá"

# New "E-mail" column header, styled like the other "special" bold header
# (same look as the "Identificador" header on the Clientes sheet).
$ws.Range("K1").Value = "E-mail"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Font.Color = 0

# E-mail value for Manjerico, wired up as a real mailto: hyperlink.
$ws.Range("K2").Value = "manjerico@manjerico.com.br"
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:manjerico@manjerico.com.br")

# Reflect the newly-added column in the sheet's selection.
[void]$ws.Range("K1:K2").Select()
